$d = $word.ActiveDocument

$pairs = @(
    @("136÷2=", "397÷8="),
    @("537÷5=", "331÷3="),
    @("154÷7=", "926÷7="),
    @("905÷3=", "472÷3="),
    @("204÷2=", "173÷9="),
    @("225÷8=", "532÷3="),
    @("110÷7=", "563÷7="),
    @("151÷4=", "683÷7="),
    @("717÷9=", "287÷8="),
    @("560÷5=", "167÷2="),
    @("222÷9=", "918÷8="),
    @("120÷8=", "243÷8="),
    @("308÷8=", "196÷5="),
    @("755÷9=", "853÷8="),
    @("201÷6=", "511÷2="),
    @("396÷3=", "888÷4="),
    @("430÷9=", "325÷4="),
    @("473÷3=", "475÷3="),
    @("257÷3=", "979÷6="),
    @("253÷4=", "212÷3="),
    @("165÷5=", "369÷8="),
    @("142÷7=", "414÷9="),
    @("835÷3=", "173÷7="),
    @("230÷8=", "881÷9="),
    @("276÷7=", "691÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
